$wb = $excel.ActiveWorkbook

# --- Meta sheet -----------------------------------------------------------
$meta = $wb.Worksheets.Item("Meta")

# The "MassBays" sampling-method-context values for TP / TDP / E.coli are
# replaced with the standard "MassWateR" context.
$meta.Range("B4").Value = "MassWateR"
$meta.Range("B5").Value = "MassWateR"
$meta.Range("B6").Value = "MassWateR"

# The data rows (A2:F6) lose their explicit cell borders; they now simply
# inherit the left-aligned / center-aligned column formatting.
$colA = $meta.Range("A2:A6")
$colA.HorizontalAlignment = -4131   # xlLeft
$colA.Borders.LineStyle = 0         # xlNone

$colBF = $meta.Range("B2:F6")
$colBF.HorizontalAlignment = -4108  # xlCenter
$colBF.Borders.LineStyle = 0        # xlNone

# Selection on the Meta sheet moves back to A2.
$meta.Range("A2").Select()

# --- Instructions sheet ----------------------------------------------------
$instr = $wb.Worksheets.Item("Instructions")

# Selection on the Instructions sheet moves to A7.
$instr.Range("A7").Select()

# Restore "Meta" as the active/visible tab (matches the saved workbook view).
$meta.Activate()
